# Add "season record" columns (Wins, Losses, Ties) to the right of the
# existing team-statistics table, per the author's commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header row (row 1): new labels in AD1:AF1 --------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the other header cells (bold, centered, bordered)
# by copying the formatting from the existing last header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# -- Data rows (2-47): season record for every player/row ---------------
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 94  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 68  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}
